$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New/changed financial-data values for rows 2-6 (columns D:AJ) ---
# Keyed by cell address -> new value.
$updates = @{
    "D2" = 7707
    "E2" = 305
    "F2" = 305
    "G2" = 241
    "H2" = 129
    "I2" = 129
    "K2" = 7945
    "L2" = 4810
    "M2" = 3135
    "N2" = 3135
    "P2" = 371
    "Q2" = 436
    "R2" = -30
    "S2" = -360
    "T2" = 130
    "U2" = 306
    "V2" = 2693
    "W2" = 3.96
    "X2" = 1.67
    "Y2" = 4.16
    "Z2" = 1.58
    "AA2" = 153.44
    "AB2" = 751.44
    "AC2" = 348
    "AD2" = 22.89
    "AE2" = 8456
    "AF2" = 0.94
    "AG2" = 80
    "AH2" = 1.01
    "AI2" = 23
    "AJ2" = 37080390
    "D3" = 7146
    "E3" = 142
    "F3" = 142
    "G3" = 96
    "H3" = 83
    "I3" = 83
    "K3" = 7354
    "L3" = 4142
    "M3" = 3213
    "N3" = 3213
    "P3" = 371
    "Q3" = 267
    "R3" = 93
    "S3" = -438
    "T3" = 45
    "U3" = 223
    "V3" = 2258
    "W3" = 1.98
    "X3" = 1.17
    "Y3" = 2.63
    "Z3" = 1.09
    "AA3" = 128.92
    "AB3" = 768.39
    "AC3" = 225
    "AD3" = 22.87
    "AE3" = 8666
    "AF3" = 0.59
    "AG3" = 80
    "AH3" = 1.56
    "AI3" = 35.58
    "AJ3" = 37080390
    "D4" = 6781
    "E4" = 45
    "F4" = 45
    "G4" = 57
    "H4" = 14
    "I4" = 14
    "K4" = 7288
    "L4" = 4088
    "M4" = 3199
    "N4" = 3199
    "P4" = 371
    "Q4" = -100
    "R4" = 163
    "S4" = -41
    "T4" = 68
    "U4" = -168
    "V4" = 2244
    "W4" = 0.66
    "X4" = 0.21
    "Y4" = 0.44
    "Z4" = 0.19
    "AA4" = 127.8
    "AB4" = 764.14
    "AC4" = 38
    "AD4" = 130.23
    "AE4" = 8629
    "AF4" = 0.58
    "AG4" = 60
    "AH4" = 1.2
    "AI4" = 156.72
    "AJ4" = 37080390
    "D5" = 7620
    "E5" = 60
    "F5" = 60
    "G5" = -18
    "H5" = -21
    "I5" = -21
    "K5" = 7568
    "L5" = 4388
    "M5" = 3180
    "N5" = 3180
    "P5" = 371
    "Q5" = 45
    "R5" = -246
    "S5" = 159
    "T5" = 149
    "U5" = -104
    "V5" = 2367
    "W5" = 0.78
    "X5" = -0.27
    "Y5" = -0.66
    "Z5" = -0.28
    "AA5" = 138
    "AB5" = 759.75
    "AC5" = -56
    "AD5" = -78.85
    "AE5" = 8577
    "AF5" = 0.52
    "AG5" = 60
    "AH5" = 1.35
    "AI5" = -106.53
    "AJ5" = 37080390
    "D6" = 7314
    "E6" = 52
    "F6" = 52
    "G6" = -131
    "H6" = -144
    "I6" = -144
    "K6" = 7786
    "L6" = 4783
    "M6" = 3003
    "N6" = 3003
    "P6" = 371
    "Q6" = -180
    "R6" = -288
    "S6" = 430
    "T6" = 215
    "U6" = -395
    "V6" = 2800
    "W6" = 0.72
    "X6" = -1.98
    "Y6" = -4.67
    "Z6" = -1.88
    "AA6" = 159.26
    "AB6" = 712.61
    "AC6" = -390
    "AD6" = -7.7
    "AE6" = 8101
    "AF6" = 0.37
    "AI6" = 0
    "AJ6" = 37080390
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Cells dropped entirely in the revised report (column realignment) ---
$clearedAddrs = @(
    "J2",
    "O2",
    "J3",
    "O3",
    "J4",
    "O4",
    "J5",
    "O5",
    "AG6",
    "AH6",
    "D7",
    "E7",
    "G7",
    "H7",
    "I7",
    "K7",
    "L7",
    "M7",
    "N7",
    "P7",
    "Q7",
    "R7",
    "S7",
    "T7",
    "U7",
    "W7",
    "X7",
    "Y7",
    "Z7",
    "AA7",
    "AC7",
    "AD7",
    "AE7",
    "AF7",
    "AG7",
    "AH7",
    "AI7",
    "D8",
    "E8",
    "G8",
    "H8",
    "I8",
    "K8",
    "L8",
    "M8",
    "N8",
    "P8",
    "Q8",
    "R8",
    "S8",
    "T8",
    "U8",
    "W8",
    "X8",
    "Y8",
    "Z8",
    "AA8",
    "AC8",
    "AD8",
    "AE8",
    "AF8",
    "AG8",
    "AH8",
    "AI8",
    "D9",
    "E9",
    "G9",
    "H9",
    "I9",
    "K9",
    "L9",
    "M9",
    "N9",
    "P9",
    "Q9",
    "R9",
    "S9",
    "T9",
    "U9",
    "W9",
    "X9",
    "Y9",
    "Z9",
    "AA9",
    "AC9",
    "AD9",
    "AE9",
    "AF9",
    "AG9",
    "AH9",
    "AI9"
)

foreach ($addr in $clearedAddrs) {
    $ws.Range($addr).ClearContents()
}
